$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Coin/Link/Price/Volume columns in this sheet are stored as plain text
# (even values that look numeric, e.g. "64.30" or "1.004"), so force the
# Text number format before writing each value to avoid Excel silently
# re-interpreting them as numbers and dropping significant trailing zeros.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.905.34'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.16%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.640.05'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.01%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.89'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.17%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5046'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.77%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.004'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.10%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2569'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.20%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06389'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.97%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.71'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.51%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07737'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.90%  '

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.653.98'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.20%  '

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.258'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.97%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.863.57'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.11%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5468'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.44%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅7897'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.57%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.30'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.12%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.910.20'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.24%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.004'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.07%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '203.04'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.81%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.376'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.99%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.896'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.969'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.005'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.07%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.926'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +9.33%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.02'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.23%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.54%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.67'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.21%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.774'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.27%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.246'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.03%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04967'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.76%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.278'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.82%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.188'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.12%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.544'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.29%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.379'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.17%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.89%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.8943'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.53%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.157.61'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.95%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5619'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.88%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01567'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.51%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.005'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.05%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.657'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.64%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.93'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.44%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8069'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.15%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.776.25'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₈118'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.85%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4544'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.10%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.002'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.44%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.86'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.17%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05058'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.57%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.002'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.26%  '
